# Append the new weekly price record (row 21) to the sheet, matching the
# existing table's layout/formatting exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

$ws.Range("A$row").Value = 8
$ws.Range("B$row").Value = "Terminal La Palmera de La Serena"
$ws.Range("C$row").Value = "Coquimbo"

$ws.Range("D$row").Value = 45013
$ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("E$row").Value = 4
$ws.Range("F$row").Value = 100112039
$ws.Range("G$row").Value = "Ciboulette"
$ws.Range("H$row").Value = "Sin especificar"
$ws.Range("I$row").Value = "Primera"
$ws.Range("J$row").Value = 1100
$ws.Range("K$row").Value = 2000
$ws.Range("L$row").Value = 2500
$ws.Range("M$row").Value = 2250
$ws.Range("N$row").Value = "$/docena de atados"
$ws.Range("O$row").Value = "Provincia del Elquí"
$ws.Range("P$row").Value = 750
$ws.Range("Q$row").Value = 3
$ws.Range("R$row").Value = "Hortaliza"
